$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 20.83749233333333
$ws.Range("H2").Value = 62.512477
$ws.Range("I2").Value = 0.02059261534406822
$ws.Range("J2").Value = 0.02059261534406822
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.502378
$ws.Range("N2").Value = 1.507134
$ws.Range("O2").Value = 0.0189742916423209
$ws.Range("P2").Value = 0.0189742916423209
$ws.Range("Q2").Value = 10.46829772343533
$ws.Range("R2").Value = 94.214679510918
$ws.Range("S2").Value = 0.0003907302892164827
$ws.Range("T2").Value = 0.0003907302892164828

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 20.83749233333333
$ws.Range("H3").Value = 62.512477
$ws.Range("I3").Value = 0.02059261534406822
$ws.Range("J3").Value = 0.02059261534406822
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6023626666666666
$ws.Range("N3").Value = 1.807088
$ws.Range("O3").Value = 0.02275060793223323
$ws.Range("P3").Value = 0.02275060793223323
$ws.Range("Q3").Value = 12.55172744855289
$ws.Range("R3").Value = 112.965547036976
$ws.Range("S3").Value = 0.0004684945179921861
$ws.Range("T3").Value = 0.0004684945179921861

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 20.83749233333333
$ws.Range("H4").Value = 62.512477
$ws.Range("I4").Value = 0.02059261534406822
$ws.Range("J4").Value = 0.02059261534406822
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.474800666666667
$ws.Range("N4").Value = 16.424402
$ws.Range("O4").Value = 0.2067774952981744
$ws.Range("P4").Value = 0.2067774952981744
$ws.Range("Q4").Value = 114.0811169181949
$ws.Range("R4").Value = 1026.730052263754
$ws.Range("S4").Value = 0.004258089422485179
$ws.Range("T4").Value = 0.00425808942248518

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.83749233333333
$ws.Range("H5").Value = 62.512477
$ws.Range("I5").Value = 0.02059261534406822
$ws.Range("J5").Value = 0.02059261534406822
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.012936333333333
$ws.Range("N5").Value = 15.038809
$ws.Range("O5").Value = 0.1893333624741797
$ws.Range("P5").Value = 0.1893333624741797
$ws.Range("Q5").Value = 104.4570224133214
$ws.Range("R5").Value = 940.1132017198931
$ws.Range("S5").Value = 0.003898869105229823
$ws.Range("T5").Value = 0.003898869105229824

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.83749233333333
$ws.Range("H6").Value = 62.512477
$ws.Range("I6").Value = 0.02059261534406822
$ws.Range("J6").Value = 0.02059261534406822
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.423250333333333
$ws.Range("N6").Value = 22.269751
$ws.Range("O6").Value = 0.2803684014001858
$ws.Range("P6").Value = 0.2803684014001858
$ws.Range("Q6").Value = 154.6819219092474
$ws.Range("R6").Value = 1392.137297183227
$ws.Range("S6").Value = 0.005773518644665342
$ws.Range("T6").Value = 0.005773518644665343

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 20.83749233333333
$ws.Range("H7").Value = 62.512477
$ws.Range("I7").Value = 0.02059261534406822
$ws.Range("J7").Value = 0.02059261534406822
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 7.461044333333334
$ws.Range("N7").Value = 22.383133
$ws.Range("O7").Value = 0.281795841252906
$ws.Range("P7").Value = 0.281795841252906
$ws.Range("Q7").Value = 155.4694540944934
$ws.Range("R7").Value = 1399.225086850441
$ws.Range("S7").Value = 0.005802913364479204
$ws.Range("T7").Value = 0.005802913364479204

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 30.85934066666667
$ws.Range("H8").Value = 92.578022
$ws.Range("I8").Value = 0.03049668942666733
$ws.Range("J8").Value = 0.03049668942666734
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.502378
$ws.Range("N8").Value = 1.507134
$ws.Range("O8").Value = 0.0189742916423209
$ws.Range("P8").Value = 0.0189742916423209
$ws.Range("Q8").Value = 15.50305384543867
$ws.Range("R8").Value = 139.527484608948
$ws.Range("S8").Value = 0.0005786530793068703
$ws.Range("T8").Value = 0.0005786530793068703

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 30.85934066666667
$ws.Range("H9").Value = 92.578022
$ws.Range("I9").Value = 0.03049668942666733
$ws.Range("J9").Value = 0.03049668942666734
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6023626666666666
$ws.Range("N9").Value = 1.807088
$ws.Range("O9").Value = 0.02275060793223323
$ws.Range("P9").Value = 0.02275060793223323
$ws.Range("Q9").Value = 18.58851473554844
$ws.Range("R9").Value = 167.296632619936
$ws.Range("S9").Value = 0.000693818224377191
$ws.Range("T9").Value = 0.0006938182243771911

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 30.85934066666667
$ws.Range("H10").Value = 92.578022
$ws.Range("I10").Value = 0.03049668942666733
$ws.Range("J10").Value = 0.03049668942666734
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.474800666666667
$ws.Range("N10").Value = 16.424402
$ws.Range("O10").Value = 0.2067774952981744
$ws.Range("P10").Value = 0.2067774952981744
$ws.Range("Q10").Value = 168.9487388547604
$ws.Range("R10").Value = 1520.538649692844
$ws.Range("S10").Value = 0.006306029054532589
$ws.Range("T10").Value = 0.00630602905453259

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 30.85934066666667
$ws.Range("H11").Value = 92.578022
$ws.Range("I11").Value = 0.03049668942666733
$ws.Range("J11").Value = 0.03049668942666734
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.012936333333333
$ws.Range("N11").Value = 15.038809
$ws.Range("O11").Value = 0.1893333624741797
$ws.Range("P11").Value = 0.1893333624741797
$ws.Range("Q11").Value = 154.6959100506442
$ws.Range("R11").Value = 1392.263190455798
$ws.Range("S11").Value = 0.00577404075348169
$ws.Range("T11").Value = 0.005774040753481691

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 30.85934066666667
$ws.Range("H12").Value = 92.578022
$ws.Range("I12").Value = 0.03049668942666733
$ws.Range("J12").Value = 0.03049668942666734
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 7.423250333333333
$ws.Range("N12").Value = 22.269751
$ws.Range("O12").Value = 0.2803684014001858
$ws.Range("P12").Value = 0.2803684014001858
$ws.Range("Q12").Value = 229.0766108902802
$ws.Range("R12").Value = 2061.689498012522
$ws.Range("S12").Value = 0.008550308062552667
$ws.Range("T12").Value = 0.00855030806255267

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 30.85934066666667
$ws.Range("H13").Value = 92.578022
$ws.Range("I13").Value = 0.03049668942666733
$ws.Range("J13").Value = 0.03049668942666734
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.461044333333334
$ws.Range("N13").Value = 22.383133
$ws.Range("O13").Value = 0.281795841252906
$ws.Range("P13").Value = 0.281795841252906
$ws.Range("Q13").Value = 230.2429088114363
$ws.Range("R13").Value = 2072.186179302926
$ws.Range("S13").Value = 0.008593840252416326
$ws.Range("T13").Value = 0.008593840252416324

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 430.5572513333333
$ws.Range("H14").Value = 1291.671754
$ws.Range("I14").Value = 0.4254974503877027
$ws.Range("J14").Value = 0.4254974503877028
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.502378
$ws.Range("N14").Value = 1.507134
$ws.Range("O14").Value = 0.0189742916423209
$ws.Range("P14").Value = 0.0189742916423209
$ws.Range("Q14").Value = 216.3024908103373
$ws.Range("R14").Value = 1946.722417293036
$ws.Range("S14").Value = 0.00807351271672024
$ws.Range("T14").Value = 0.008073512716720242

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 430.5572513333333
$ws.Range("H15").Value = 1291.671754
$ws.Range("I15").Value = 0.4254974503877027
$ws.Range("J15").Value = 0.4254974503877028
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.6023626666666666
$ws.Range("N15").Value = 1.807088
$ws.Range("O15").Value = 0.02275060793223323
$ws.Range("P15").Value = 0.02275060793223323
$ws.Range("Q15").Value = 259.3516140658169
$ws.Range("R15").Value = 2334.164526592352
$ws.Range("S15").Value = 0.009680325669935484
$ws.Range("T15").Value = 0.009680325669935484

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 430.5572513333333
$ws.Range("H16").Value = 1291.671754
$ws.Range("I16").Value = 0.4254974503877027
$ws.Range("J16").Value = 0.4254974503877028
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.474800666666667
$ws.Range("N16").Value = 16.424402
$ws.Range("O16").Value = 0.2067774952981744
$ws.Range("P16").Value = 0.2067774952981744
$ws.Range("Q16").Value = 2357.215126637901
$ws.Range("R16").Value = 21214.93613974111
$ws.Range("S16").Value = 0.08798329704692838
$ws.Range("T16").Value = 0.0879832970469284

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 430.5572513333333
$ws.Range("H17").Value = 1291.671754
$ws.Range("I17").Value = 0.4254974503877027
$ws.Range("J17").Value = 0.4254974503877028
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.012936333333333
$ws.Range("N17").Value = 15.038809
$ws.Range("O17").Value = 0.1893333624741797
$ws.Range("P17").Value = 0.1893333624741797
$ws.Range("Q17").Value = 2158.356088788998
$ws.Range("R17").Value = 19425.20479910099
$ws.Range("S17").Value = 0.08056086300609422
$ws.Range("T17").Value = 0.08056086300609423

$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 430.5572513333333
$ws.Range("H18").Value = 1291.671754
$ws.Range("I18").Value = 0.4254974503877027
$ws.Range("J18").Value = 0.4254974503877028
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 7.423250333333333
$ws.Range("N18").Value = 22.269751
$ws.Range("O18").Value = 0.2803684014001858
$ws.Range("P18").Value = 0.2803684014001858
$ws.Range("Q18").Value = 3196.13425947925
$ws.Range("R18").Value = 28765.20833531325
$ws.Range("S18").Value = 0.1192960399650551
$ws.Range("T18").Value = 0.1192960399650551

$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 430.5572513333333
$ws.Range("H19").Value = 1291.671754
$ws.Range("I19").Value = 0.4254974503877027
$ws.Range("J19").Value = 0.4254974503877028
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 7.461044333333334
$ws.Range("N19").Value = 22.383133
$ws.Range("O19").Value = 0.281795841252906
$ws.Range("P19").Value = 0.281795841252906
$ws.Range("Q19").Value = 3212.406740236142
$ws.Range("R19").Value = 28911.66066212528
$ws.Range("S19").Value = 0.1199034119829693
$ws.Range("T19").Value = 0.1199034119829693

$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 371.6201173333334
$ws.Range("H20").Value = 1114.860352
$ws.Range("I20").Value = 0.3672529308203304
$ws.Range("J20").Value = 0.3672529308203304
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 0.502378
$ws.Range("N20").Value = 1.507134
$ws.Range("O20").Value = 0.0189742916423209
$ws.Range("P20").Value = 0.0189742916423209
$ws.Range("Q20").Value = 186.6937713056853
$ws.Range("R20").Value = 1680.243941751168
$ws.Range("S20").Value = 0.006968364215882051
$ws.Range("T20").Value = 0.006968364215882052

$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 371.6201173333334
$ws.Range("H21").Value = 1114.860352
$ws.Range("I21").Value = 0.3672529308203304
$ws.Range("J21").Value = 0.3672529308203304
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 0.6023626666666666
$ws.Range("N21").Value = 1.807088
$ws.Range("O21").Value = 0.02275060793223323
$ws.Range("P21").Value = 0.02275060793223323
$ws.Range("Q21").Value = 223.8500848638862
$ws.Range("R21").Value = 2014.650763774976
$ws.Range("S21").Value = 0.00835522744105691
$ws.Range("T21").Value = 0.008355227441056911

$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 371.6201173333334
$ws.Range("H22").Value = 1114.860352
$ws.Range("I22").Value = 0.3672529308203304
$ws.Range("J22").Value = 0.3672529308203304
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 5.474800666666667
$ws.Range("N22").Value = 16.424402
$ws.Range("O22").Value = 0.2067774952981744
$ws.Range("P22").Value = 0.2067774952981744
$ws.Range("Q22").Value = 2034.546066123278
$ws.Range("R22").Value = 18310.91459510951
$ws.Range("S22").Value = 0.07593964117594162
$ws.Range("T22").Value = 0.07593964117594164

$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 371.6201173333334
$ws.Range("H23").Value = 1114.860352
$ws.Range("I23").Value = 0.3672529308203304
$ws.Range("J23").Value = 0.3672529308203304
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 5.012936333333333
$ws.Range("N23").Value = 15.038809
$ws.Range("O23").Value = 0.1893333624741797
$ws.Range("P23").Value = 0.1893333624741797
$ws.Range("Q23").Value = 1862.907988377863
$ws.Range("R23").Value = 16766.17189540077
$ws.Range("S23").Value = 0.06953323227071047
$ws.Range("T23").Value = 0.06953323227071048

$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 371.6201173333334
$ws.Range("H24").Value = 1114.860352
$ws.Range("I24").Value = 0.3672529308203304
$ws.Range("J24").Value = 0.3672529308203304
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 7.423250333333333
$ws.Range("N24").Value = 22.269751
$ws.Range("O24").Value = 0.2803684014001858
$ws.Range("P24").Value = 0.2803684014001858
$ws.Range("Q24").Value = 2758.629159868039
$ws.Range("R24").Value = 24827.66243881236
$ws.Range("S24").Value = 0.102966117123629
$ws.Range("T24").Value = 0.1029661171236291

$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 371.6201173333334
$ws.Range("H25").Value = 1114.860352
$ws.Range("I25").Value = 0.3672529308203304
$ws.Range("J25").Value = 0.3672529308203304
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 7.461044333333334
$ws.Range("N25").Value = 22.383133
$ws.Range("O25").Value = 0.281795841252906
$ws.Range("P25").Value = 0.281795841252906
$ws.Range("Q25").Value = 2772.674170582536
$ws.Range("R25").Value = 24954.06753524282
$ws.Range("S25").Value = 0.1034903485931103
$ws.Range("T25").Value = 0.1034903485931103

$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 145.447436
$ws.Range("H26").Value = 436.342308
$ws.Range("I26").Value = 0.1437381741726046
$ws.Range("J26").Value = 0.1437381741726046
$ws.Range("K26").Value = 3
$ws.Range("M26").Value = 0.502378
$ws.Range("N26").Value = 1.507134
$ws.Range("O26").Value = 0.0189742916423209
$ws.Range("P26").Value = 0.0189742916423209
$ws.Range("Q26").Value = 73.06959200280801
$ws.Range("R26").Value = 657.626328025272
$ws.Range("S26").Value = 0.002727330036885718
$ws.Range("T26").Value = 0.002727330036885718

$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 145.447436
$ws.Range("H27").Value = 436.342308
$ws.Range("I27").Value = 0.1437381741726046
$ws.Range("J27").Value = 0.1437381741726046
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 0.6023626666666666
$ws.Range("N27").Value = 1.807088
$ws.Range("O27").Value = 0.02275060793223323
$ws.Range("P27").Value = 0.02275060793223323
$ws.Range("Q27").Value = 87.61210540878933
$ws.Range("R27").Value = 788.5089486791039
$ws.Range("S27").Value = 0.00327013084549598
$ws.Range("T27").Value = 0.00327013084549598

$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 145.447436
$ws.Range("H28").Value = 436.342308
$ws.Range("I28").Value = 0.1437381741726046
$ws.Range("J28").Value = 0.1437381741726046
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 5.474800666666667
$ws.Range("N28").Value = 16.424402
$ws.Range("O28").Value = 0.2067774952981744
$ws.Range("P28").Value = 0.2067774952981744
$ws.Range("Q28").Value = 796.2957195777574
$ws.Range("R28").Value = 7166.661476199816
$ws.Range("S28").Value = 0.02972181963414392
$ws.Range("T28").Value = 0.02972181963414392

$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 145.447436
$ws.Range("H29").Value = 436.342308
$ws.Range("I29").Value = 0.1437381741726046
$ws.Range("J29").Value = 0.1437381741726046
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 5.012936333333333
$ws.Range("N29").Value = 15.038809
$ws.Range("O29").Value = 0.1893333624741797
$ws.Range("P29").Value = 0.1893333624741797
$ws.Range("Q29").Value = 729.1187365145747
$ws.Range("R29").Value = 6562.068628631172
$ws.Range("S29").Value = 0.02721443183199853
$ws.Range("T29").Value = 0.02721443183199853

$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 145.447436
$ws.Range("H30").Value = 436.342308
$ws.Range("I30").Value = 0.1437381741726046
$ws.Range("J30").Value = 0.1437381741726046
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 7.423250333333333
$ws.Range("N30").Value = 22.269751
$ws.Range("O30").Value = 0.2803684014001858
$ws.Range("P30").Value = 0.2803684014001858
$ws.Range("Q30").Value = 1079.692727769479
$ws.Range("R30").Value = 9717.234549925308
$ws.Range("S30").Value = 0.04029964211295463
$ws.Range("T30").Value = 0.04029964211295463

$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 145.447436
$ws.Range("H31").Value = 436.342308
$ws.Range("I31").Value = 0.1437381741726046
$ws.Range("J31").Value = 0.1437381741726046
$ws.Range("K31").Value = 3
$ws.Range("M31").Value = 7.461044333333334
$ws.Range("N31").Value = 22.383133
$ws.Range("O31").Value = 0.281795841252906
$ws.Range("P31").Value = 0.281795841252906
$ws.Range("Q31").Value = 1085.189768165663
$ws.Range("R31").Value = 9766.707913490964
$ws.Range("S31").Value = 0.04050481971112585
$ws.Range("T31").Value = 0.04050481971112584

$ws.Range("E32").Value = 3
$ws.Range("G32").Value = 12.56985766666667
$ws.Range("H32").Value = 37.709573
$ws.Range("I32").Value = 0.01242213984862671
$ws.Range("J32").Value = 0.01242213984862671
$ws.Range("K32").Value = 3
$ws.Range("M32").Value = 0.502378
$ws.Range("N32").Value = 1.507134
$ws.Range("O32").Value = 0.0189742916423209
$ws.Range("P32").Value = 0.0189742916423209
$ws.Range("Q32").Value = 6.314819954864666
$ws.Range("R32").Value = 56.833379593782
$ws.Range("S32").Value = 0.0002357013043095391
$ws.Range("T32").Value = 0.0002357013043095392

$ws.Range("E33").Value = 3
$ws.Range("G33").Value = 12.56985766666667
$ws.Range("H33").Value = 37.709573
$ws.Range("I33").Value = 0.01242213984862671
$ws.Range("J33").Value = 0.01242213984862671
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 0.6023626666666666
$ws.Range("N33").Value = 1.807088
$ws.Range("O33").Value = 0.02275060793223323
$ws.Range("P33").Value = 0.02275060793223323
$ws.Range("Q33").Value = 7.571612983713776
$ws.Range("R33").Value = 68.144516853424
$ws.Range("S33").Value = 0.0002826112333754772
$ws.Range("T33").Value = 0.0002826112333754772

$ws.Range("E34").Value = 3
$ws.Range("G34").Value = 12.56985766666667
$ws.Range("H34").Value = 37.709573
$ws.Range("I34").Value = 0.01242213984862671
$ws.Range("J34").Value = 0.01242213984862671
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 5.474800666666667
$ws.Range("N34").Value = 16.424402
$ws.Range("O34").Value = 0.2067774952981744
$ws.Range("P34").Value = 0.2067774952981744
$ws.Range("Q34").Value = 68.81746513337177
$ws.Range("R34").Value = 619.357186200346
$ws.Range("S34").Value = 0.002568618964142673
$ws.Range("T34").Value = 0.002568618964142674

$ws.Range("E35").Value = 3
$ws.Range("G35").Value = 12.56985766666667
$ws.Range("H35").Value = 37.709573
$ws.Range("I35").Value = 0.01242213984862671
$ws.Range("J35").Value = 0.01242213984862671
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 5.012936333333333
$ws.Range("N35").Value = 15.038809
$ws.Range("O35").Value = 0.1893333624741797
$ws.Range("P35").Value = 0.1893333624741797
$ws.Range("Q35").Value = 63.01189620206188
$ws.Range("R35").Value = 567.107065818557
$ws.Range("S35").Value = 0.002351925506664992
$ws.Range("T35").Value = 0.002351925506664992

$ws.Range("E36").Value = 3
$ws.Range("G36").Value = 12.56985766666667
$ws.Range("H36").Value = 37.709573
$ws.Range("I36").Value = 0.01242213984862671
$ws.Range("J36").Value = 0.01242213984862671
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 7.423250333333333
$ws.Range("N36").Value = 22.269751
$ws.Range("O36").Value = 0.2803684014001858
$ws.Range("P36").Value = 0.2803684014001858
$ws.Range("Q36").Value = 93.30920011403587
$ws.Range("R36").Value = 839.782801026323
$ws.Range("S36").Value = 0.003482775491329015
$ws.Range("T36").Value = 0.003482775491329016

$ws.Range("E37").Value = 3
$ws.Range("G37").Value = 12.56985766666667
$ws.Range("H37").Value = 37.709573
$ws.Range("I37").Value = 0.01242213984862671
$ws.Range("J37").Value = 0.01242213984862671
$ws.Range("K37").Value = 3
$ws.Range("M37").Value = 7.461044333333334
$ws.Range("N37").Value = 22.383133
$ws.Range("O37").Value = 0.281795841252906
$ws.Range("P37").Value = 0.281795841252906
$ws.Range("Q37").Value = 93.78426531468989
$ws.Range("R37").Value = 844.058387832209
$ws.Range("S37").Value = 0.003500507348805009
$ws.Range("T37").Value = 0.003500507348805009
